$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows for "Huelva" and "Huesca" swap places: row 53 (previously Huelva,
# Casos activos=72) becomes "Huesca" with Casos activos=0, and row 54
# (previously Huesca, Casos activos=0) becomes "Huelva" with Casos activos=72.
$ws.Range("A53").Value = "Huesca"
$ws.Range("C53").Value = 0

$ws.Range("A54").Value = "Huelva"
$ws.Range("C54").Value = 72

# Update the "last updated" timestamp from 02:16 to 02:46.
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 02:46"
